$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.237.87"
$ws.Range("E2").Value = "  +1.08%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.571.28"
$ws.Range("E3").Value = "  +0.72%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.28%  "

# Row 5 - BNB (numeric-looking price -> force text with quote prefix)
$ws.Range("D5").Value = "'210.97"
$ws.Range("E5").Value = "  +1.86%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.492"
$ws.Range("E6").Value = "  +0.61%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.20%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'22.04"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.249"
$ws.Range("E9").Value = "  +0.36%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.04%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0868"
$ws.Range("E11").Value = "  +1.24%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.794.46"
$ws.Range("E12").Value = "  +0.66%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.552.14"
$ws.Range("E13").Value = "  -0.57%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'3.79"
$ws.Range("E14").Value = "  +0.84%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.06%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "27.199.76"
$ws.Range("E16").Value = "  +0.85%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'62.24"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'7.55"
$ws.Range("E18").Value = "  +2.75%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'216.42"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.42%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.29%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  +1.22%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "'9.23"
$ws.Range("E23").Value = "  +0.34%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.54%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'153.87"
$ws.Range("E25").Value = "  +0.68%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'6.64"
$ws.Range("E26").Value = "  +0.46%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'15.08"
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +2.24%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.25%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +2.63%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +0.46%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.17%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.450.48"
$ws.Range("E33").Value = "  +2.21%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +1.41%  "

# Row 35 - TrustWalletToken
$ws.Range("E35").Value = "  +7.03%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +0.38%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.55%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.0167"
$ws.Range("E38").Value = "  +0.94%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "'0.534"
$ws.Range("E39").Value = "  +0.43%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "'5.87"
$ws.Range("E40").Value = "  +2.51%  "

# Row 41 - ARBITRUM
$ws.Range("D41").Value = "'0.811"
$ws.Range("E41").Value = "  +0.35%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.28%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +1.00%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value = "  -0.62%  "

# Row 45 - Aave
$ws.Range("D45").Value = "'64.42"
$ws.Range("E45").Value = "  -0.43%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  -1.27%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.705.63"
$ws.Range("E47").Value = "  +0.53%  "

# Row 48 - Quant
$ws.Range("D48").Value = "'86.04"
$ws.Range("E48").Value = "  -1.59%  "

# Row 49/50 swap: Cronos moves to row 49, BabyDogeCoin moves to row 50,
# each with updated price/volume values.
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0525"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  +1.09%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -0.26%  "
